# Updated cryptos list on Fri Feb  2 03:31:13 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force Excel to store the value as literal text (matches the source
    # feed's inline strings) instead of silently re-interpreting
    # numeric-looking strings like "302.24" or "49.10" as numbers, which
    # would corrupt formatting (trailing zeros, thousand-dot grouping, etc).
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "43.062.94"
$ws.Range("E2").Value = "  +2.50%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.302.78"
$ws.Range("E3").Value = "  +2.03%  "

# Row 5 - BNB
Set-TextCell "D5" "302.24"
$ws.Range("E5").Value = "  +1.37%  "

# Row 6 - Solana
Set-TextCell "D6" "99.55"
$ws.Range("E6").Value = "  +6.35%  "

# Row 7 - XRP
Set-TextCell "D7" "0.506"
$ws.Range("E7").Value = "  +1.64%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +3.29%  "

# Row 10 - Avalanche
Set-TextCell "D10" "34.42"

# Row 11 - Dogecoin
Set-TextCell "D11" "0.0799"
$ws.Range("E11").Value = "  +1.60%  "

# Row 12 - OKB
Set-TextCell "D12" "49.10"
$ws.Range("E12").Value = "  +3.92%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +4.15%  "

# Row 14 - Chainlink
Set-TextCell "D14" "17.71"
$ws.Range("E14").Value = "  +16.09%  "

# Row 15 - Polkadot
Set-TextCell "D15" "6.80"
$ws.Range("E15").Value = "  +2.41%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextCell "D16" "2.662.23"
$ws.Range("E16").Value = "  +2.00%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.326.53"
$ws.Range("E17").Value = "  +2.91%  "

# Row 18 - Polygon
Set-TextCell "D18" "0.808"
$ws.Range("E18").Value = "  +4.47%  "

# Row 19 - WrappedBTC
Set-TextCell "D19" "42.956.61"
$ws.Range("E19").Value = "  +2.17%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextCell "D20" "12.32"
$ws.Range("E20").Value = "  +8.52%  "

# Row 21 - ShibaInu
Set-TextCell "D21" "0.0₃0907"
$ws.Range("E21").Value = "  +2.05%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.60%  "

# Row 23 - Litecoin
Set-TextCell "D23" "67.86"

# Row 24 - BitcoinCash
Set-TextCell "D24" "237.04"
$ws.Range("E24").Value = "  +1.77%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  +14.44%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +0.64%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.12%  "

# Row 28 - EthereumClassic
Set-TextCell "D28" "24.55"
$ws.Range("E28").Value = "  +3.65%  "

# Row 29 - Monero
Set-TextCell "D29" "168.76"
$ws.Range("E29").Value = "  +1.35%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -3.07%  "

# Row 31 - InjectiveProtocol
Set-TextCell "D31" "33.89"
$ws.Range("E31").Value = "  +1.22%  "

# Row 32 - Cosmos
Set-TextCell "D32" "9.19"
$ws.Range("E32").Value = "  +1.80%  "

# Row 33 - FirstDigitalUSD
Set-TextCell "D33" "0.999"
$ws.Range("E33").Value = "  -0.06%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +1.96%  "

# Row 35 - WEMIXToken
Set-TextCell "D35" "2.44"
$ws.Range("E35").Value = "  +3.88%  "

# Row 36 - RenderToken
Set-TextCell "D36" "4.56"
$ws.Range("E36").Value = "  +4.58%  "

# Row 37 - Celestia
Set-TextCell "D37" "17.09"
$ws.Range("E37").Value = "  +7.78%  "

# Row 38 - Hedera
Set-TextCell "D38" "0.0700"
$ws.Range("E38").Value = "  +1.11%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +3.70%  "

# Rows 40 & 41 swap places: LidoDAOToken <-> ARBITRUM
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D40" "1.79"
$ws.Range("E40").Value = "  +5.04%  "

$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D41" "2.82"
$ws.Range("E41").Value = "  +0.96%  "

# Row 42 - Stellar
Set-TextCell "D42" "0.109"
$ws.Range("E42").Value = "  +0.31%  "

# Row 43 - ApeXProtocol
Set-TextCell "D43" "2.35"
$ws.Range("E43").Value = "  -2.27%  "

# Row 44 - Maker
Set-TextCell "D44" "2.001.26"
$ws.Range("E44").Value = "  +2.96%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +2.66%  "

# Row 46 - FraxShare
Set-TextCell "D46" "10.13"
$ws.Range("E46").Value = "  +6.31%  "

# Row 47 - EnergySwap
Set-TextCell "D47" "17.77"
$ws.Range("E47").Value = "  +2.12%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  +3.18%  "

# Row 49 - MultiversX
Set-TextCell "D49" "55.68"
$ws.Range("E49").Value = "  +6.84%  "

# Row 50 - RocketPoolETH
Set-TextCell "D50" "2.528.46"
$ws.Range("E50").Value = "  +1.78%  "

# Row 51 - Stacks
Set-TextCell "D51" "1.52"
$ws.Range("E51").Value = "  +2.20%  "
